# Replace all occurrences of the old Ref Des "GA05MOAS-GL002" with the
# corrected "GA05MOAS-GL495" across the workbook (Moorings + Asset_Cal_Info
# sheets), preserving the suffixes used for each sub-asset row.

$wb = $excel.ActiveWorkbook

$oldPrefix = "GA05MOAS-GL002"
$newPrefix = "GA05MOAS-GL495"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldPrefix)) {
            $cell.Value = $val.Replace($oldPrefix, $newPrefix)
        }
    }
}
